$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Figure 1 table (rows 5-6): replace numeric win-rates with
#    descriptive percentage strings. Call order controls the order
#    new shared strings are appended, so keep this exact sequence.
# ------------------------------------------------------------------
$ws.Range("C6").Value = "37.86% (53.57%)"
$ws.Range("B5").Value = "38.57% (56%, random 20%, random 33% if LLM is smart)"
$ws.Range("C5").Value = "49.29% (47.86%)"
$ws.Range("B6").Value = "47.06% (57.35%)"

# ------------------------------------------------------------------
# 2) Figure 2 table: update the "2 to 4" summary rows (15-18)
# ------------------------------------------------------------------
$ws.Range("C15").Value = 0.3143
$ws.Range("C16").Value = 0.5143
$ws.Range("D16").Value = 0.0857
$ws.Range("A17").Value = 32
$ws.Range("C17").Value = 0.6875
$ws.Range("D17").Value = 0.125
$ws.Range("A18").Value = 9
$ws.Range("C18").Value = 0.7778
$ws.Range("D18").Value = 0.2222

# Rows 19-24 no longer hold game data - clear the values but keep formatting
$ws.Range("A19:D24").ClearContents()

# ------------------------------------------------------------------
# 3) New ablation study tables (rows 29-41 and 43-55)
#    Build them by cloning the formatting of the existing Figure 2
#    table (rows 14-25) and then filling in the new values.
# ------------------------------------------------------------------
$ws.Range("A14:H25").Copy()
$ws.Range("A30").PasteSpecial(-4122) | Out-Null
$ws.Range("A14:H25").Copy()
$ws.Range("A44").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Title rows 29 and 43 use the same style as row 9 / row 1 section headers
$ws.Range("A9").Copy()
$ws.Range("A29").PasteSpecial(-4122) | Out-Null
$ws.Range("A9").Copy()
$ws.Range("A43").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A29").Value = "Ablation Study: Only Voting Information"
$ws.Range("A43").Value = "Ablation Study: Without Voting Information"

# --- Table "Only Voting Information" (rows 30-41) ---
$ws.Range("A31").Value = 35
$ws.Range("B31").Value = 1
$ws.Range("C31").Value = 0.6286
$ws.Range("D31").Value = 0

$ws.Range("A32").Value = 35
$ws.Range("B32").Value = 2
$ws.Range("C32").Value = 0.4
$ws.Range("D32").Value = 0.0571
$ws.Range("E32").Value = 0.0667
$ws.Range("F32").Value = 0
$ws.Range("G32").Value = 0.3173
$ws.Range("H32").Value = 0.0153

$ws.Range("A33").Value = 32
$ws.Range("B33").Value = 3
$ws.Range("C33").Value = 0.8438
$ws.Range("D33").Value = 0.1875
$ws.Range("E33").Value = 0.3333
$ws.Range("F33").Value = 0
$ws.Range("G33").Value = 0.3267
$ws.Range("H33").Value = 0.0191

$ws.Range("A34").Value = 9
$ws.Range("B34").Value = 4
$ws.Range("C34").Value = 0.7778
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0.5
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 0.3875
$ws.Range("H34").Value = 0.02

# fix style on E32:E34 to match the new table (solid-fill percent style, no border change)
$ws.Range("G32").Copy()
$ws.Range("E32:E34").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("E32").Value = 0.0667
$ws.Range("E33").Value = 0.3333
$ws.Range("E34").Value = 0.5

$ws.Range("B41").Value = "2 to 4"
$ws.Range("C41").Formula = "=AVERAGE(C32:C34)"
$ws.Range("D41").Formula = "=AVERAGE(D32:D34)"
$ws.Range("E41").Formula = "=AVERAGE(E32:E34)"
$ws.Range("F41").Formula = "=AVERAGE(F32:F34)"
$ws.Range("G41").Formula = "=AVERAGE(G32:G34)"
$ws.Range("H41").Formula = "=AVERAGE(H32:H34)"

# --- Table "Without Voting Information" (rows 44-55) ---
$ws.Range("A45").Value = 35
$ws.Range("B45").Value = 1
$ws.Range("C45").Value = 0.3714
$ws.Range("D45").Value = 0.0286

$ws.Range("A46").Value = 35
$ws.Range("B46").Value = 2
$ws.Range("C46").Value = 0.4286
$ws.Range("D46").Value = 0.0857
$ws.Range("E46").Value = 0.4667
$ws.Range("F46").Value = 0.0667
$ws.Range("G46").Value = 0.3173
$ws.Range("H46").Value = 0.0153

$ws.Range("A47").Value = 32
$ws.Range("B47").Value = 3
$ws.Range("C47").Value = 0.7188
$ws.Range("D47").Value = 0.125
$ws.Range("E47").Value = 0.8333
$ws.Range("F47").Value = 0.1667
$ws.Range("G47").Value = 0.3267
$ws.Range("H47").Value = 0.0191

$ws.Range("A48").Value = 9
$ws.Range("B48").Value = 4
$ws.Range("C48").Value = 0.8889
$ws.Range("D48").Value = 0.3333
$ws.Range("E48").Value = 0.75
$ws.Range("F48").Value = 0
$ws.Range("G48").Value = 0.3875
$ws.Range("H48").Value = 0.02

$ws.Range("G46").Copy()
$ws.Range("E46:E48").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("E46").Value = 0.4667
$ws.Range("E47").Value = 0.8333
$ws.Range("E48").Value = 0.75

$ws.Range("B55").Value = "2 to 4"
$ws.Range("C55").Formula = "=AVERAGE(C46:C48)"
$ws.Range("D55").Formula = "=AVERAGE(D46:D48)"
$ws.Range("E55").Formula = "=AVERAGE(E46:E48)"
$ws.Range("F55").Formula = "=AVERAGE(F46:F48)"
$ws.Range("G55").Formula = "=AVERAGE(G46:G48)"
$ws.Range("H55").Formula = "=AVERAGE(H46:H48)"

# ------------------------------------------------------------------
# 4) Cosmetic workbook/view changes
# ------------------------------------------------------------------
$ws.Range("B:B").ColumnWidth = 49
$ws.Range("F16").Select()

Write-Host "Edit applied"
